$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Simple single-value cell updates
$t.Cell(1, 1).Range.Text = "0M"
$t.Cell(2, 1).Range.Text = "0M"
$t.Cell(3, 1).Range.Text = "0M"
$t.Cell(4, 1).Range.Text = "512"
$t.Cell(5, 1).Range.Text = "0.00002"
$t.Cell(6, 1).Range.Text = "0.00006"
$t.Cell(7, 1).Range.Text = "0.00004"
$t.Cell(9, 1).Range.Text = "0.00005"
$t.Cell(10, 1).Range.Text = "0.00005"
$t.Cell(11, 1).Range.Text = "0.00005"
$t.Cell(12, 1).Range.Text = "0.01769"

# Collapse the multi-run tab-separated rows down to their first value
$t.Cell(44, 1).Range.Text = "100"
$t.Cell(45, 1).Range.Text = "0.02"
$t.Cell(46, 1).Range.Text = "362"
